$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update fuel consumption (litres) and cost (kDKK) for 2021-06-19 row
$ws.Range("C4").Value = 1381
$ws.Range("D4").Value = 14

# Update total row fuel consumption (litres) and cost (kDKK)
$ws.Range("C9").Value = 13525
$ws.Range("D9").Value = 137.2
